$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next quarterly data point (01-04-2021) as row 71.
# Column A holds the period label as text ("Serie"); enter it via a
# formula that yields a text result, then convert it to a literal value
# with PasteSpecial so Excel doesn't re-interpret the text as a date.
$ws.Range("A71").Formula = "=""01-04-2021"""
$ws.Range("A71").Copy()
$ws.Range("A71").PasteSpecial(-4163)

$ws.Range("B71").Value = 2395184
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 53306
$ws.Range("F71").Value = 2323987
$ws.Range("G71").Value = 17891
